# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-6, columns B..G
$data = @{
    2 = @(0.6545652718822623, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 0, 2.964545797025059)
    3 = @(0.01253208636536152, 0.3048912486333797, 0.1496068669990043, 13.86384647080068, 1, 14.33087667279843)
    4 = @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 3.755628166162433)
    5 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 0, 6.15379541431027)
    6 = @(1.445647641019636, 1.626987699542094, 0.7210945179870265, 13.86384647080068, 1, 17.65757632934944)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 2 + $i  # Column B = 2
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
